$d = $word.ActiveDocument

function Replace-ExactText($findText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Select()
        $word.Selection.TypeText($newText)
    }
    return $found
}

Replace-ExactText " el carrete viral transmitió correctamente el atractivo del producto a través de objetos visuales atractivos e informativos." " el `"reel`" viral transmitió correctamente el atractivo del producto a través de imágenes visuales atractivas y contenido informativo."

Replace-ExactText "Marketing de influencia:" "Marketing de`"influencers`":"

Replace-ExactText " el poder del marketing de influenciadores no puede ser sobreestado." " el poder del marketing de los `"influencers`" no se puede sobrevalorar."

Replace-ExactText "Variedades de sabor y sabor:" "Variedad de gustos y sabores:"

Replace-ExactText " la reputación de Contoso Protein Plus para deliciosos y diversos sabores fue un punto de venta clave en el contenido viral." " la reputación de Contoso Protein Plus por sus sabores deliciosos y variados fue un punto de venta clave en el contenido viral."

Replace-ExactText " el aumento continuo de la conciencia de salud y fitness, combinado con un aumento en el número de personas que adoptan rutinas de entrenamiento y estilos de vida activos, crearon un mercado receptivo para un producto como Contoso Protein Plus." " el continuo aumento del interés por la salud y el ejercicio, combinado con un aumento en el número de personas que adoptan rutinas de entrenamiento y estilos de vida activos, crearon un mercado receptivo para un producto como Contoso Protein Plus."

Replace-ExactText " la accesibilidad del producto a través de varios minoristas en línea ha alimentado aún más la hype." " la accesibilidad del producto a través de varios minoristas en línea ha impulsado aún más el revuelo."

Replace-ExactText "Opiniones positivas y testimonios:" "Reseñas y testimonios positivos:"

Replace-ExactText " El carrete no era un caso aislado." " el `"reel`" no fue un caso aislado."

Replace-ExactText "Palabra de boca:" "El boca a boca:"

Replace-ExactText " las plataformas de medios sociales fomentan la rápida propagación de tendencias a través de palabras de boca." " las plataformas de redes sociales fomentan la rápida propagación de tendencias a través del boca a boca."
